# Corrected some things; removed kcal data; removed historic data for german plots
#
# The table had a row for "EU28+CH" (non-German) geographic scope mixed in
# with the German (DE) rows, and four rows reporting meat consumption in
# kcal/cap/day (historic-style duplicate units), which are removed. One
# "g/cap/day" meat-consumption row is kept.
#
# Deleting bottom-to-top (using the ORIGINAL row numbers) so that earlier
# deletions don't shift the row numbers of rows still to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "kcal meat/cap/day" row for UBA (2020) / Agriculture
# (row 29) - keep the "g/cap/day" row for the same source (row 28).
$ws.Rows(29).Delete()

# Remove the "kcal meat/cap/day" row for négaWatt et al. (2023) / Food (row 27)
$ws.Rows(27).Delete()

# Remove the "kcal/cap/day" row for Costa et al. (2021) / EUCalc / Food (row 26)
$ws.Rows(26).Delete()

# Remove the "kcal/cap/day" row for Climact (no year) / Food (row 25)
$ws.Rows(25).Delete()

# Remove the EU28+CH (non-German) duplicate of the Costa et al. energy row (row 3)
$ws.Rows(3).Delete()
